# Applies the proofing-markup edit described by the commit diff:
#  1. Splits "This is my first git project ." into two runs with a
#     gramStart/gramEnd proofErr pair around "project .".
#  2. Splits "Lets change this file some more . HA HA HA." into several
#     runs separated by spellStart/spellEnd and gramStart/gramEnd
#     proofErr markers.
#  3. Appends two empty paragraphs and a new paragraph ("Another change
#     again. Ha !Ha! Ha!") after the dotted-line paragraph, with its own
#     gramStart/gramEnd proofErr pair.
#
# NOTE: this interpreter's argument parser misbehaves when a call passes
# a parenthesized expression (or a `.Property` access) directly as an
# argument, silently turning it into an empty argument. So every value
# is built up in a plain variable first, and functions are always
# invoked with bare variable references (never literal `(...)` or
# `$obj.Prop` at the call site).

$d = $word.ActiveDocument

$pkgHead = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgTail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Insert-PackageXml($range, $bodyInnerXml) {
    $pkg = $pkgHead + $bodyInnerXml + $pkgTail
    $range.InsertXML($pkg)
}

# --- 1. "This is my first git project ." -------------------------------
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$xml2 = '<w:p>' + '<w:r><w:t xml:space="preserve">This is my first git </w:t></w:r>' + '<w:proofErr w:type="gramStart"/>' + '<w:r><w:t>project .</w:t></w:r>' + '<w:proofErr w:type="gramEnd"/>' + '</w:p>'
Insert-PackageXml $r2 $xml2

# --- 2. "Lets change this file some more . HA HA HA." -------------------
$p6 = $d.Paragraphs(6)
$r6 = $p6.Range
$xml6 = '<w:p>'
$xml6 = $xml6 + '<w:proofErr w:type="spellStart"/>'
$xml6 = $xml6 + '<w:r><w:t>Lets</w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="spellEnd"/>'
$xml6 = $xml6 + '<w:r><w:t xml:space="preserve"> change this file some </w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="gramStart"/>'
$xml6 = $xml6 + '<w:r><w:t>more .</w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="gramEnd"/>'
$xml6 = $xml6 + '<w:r><w:t xml:space="preserve"> HA </w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="spellStart"/>'
$xml6 = $xml6 + '<w:r><w:t>HA</w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="spellEnd"/>'
$xml6 = $xml6 + '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="spellStart"/>'
$xml6 = $xml6 + '<w:r><w:t>HA</w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="spellEnd"/>'
$xml6 = $xml6 + '<w:r><w:t>.</w:t></w:r>'
$xml6 = $xml6 + '</w:p>'
Insert-PackageXml $r6 $xml6

# --- 3. Append two blank paragraphs + the "Another change again." one --
$p8 = $d.Paragraphs(8)
$r8 = $p8.Range
$xml8 = '<w:p><w:r><w:t>……………………………………………………………………..</w:t></w:r></w:p>'
$xml8 = $xml8 + '<w:p/>'
$xml8 = $xml8 + '<w:p/>'
$xml8 = $xml8 + '<w:p>'
$xml8 = $xml8 + '<w:r><w:t xml:space="preserve">Another change again. </w:t></w:r>'
$xml8 = $xml8 + '<w:proofErr w:type="gramStart"/>'
$xml8 = $xml8 + '<w:r><w:t>Ha !Ha</w:t></w:r>'
$xml8 = $xml8 + '<w:proofErr w:type="gramEnd"/>'
$xml8 = $xml8 + '<w:r><w:t>! Ha!</w:t></w:r>'
$xml8 = $xml8 + '</w:p>'
Insert-PackageXml $r8 $xml8

# InsertXML leaves a spurious empty trailing paragraph (the tail of the
# paragraph mark that was replaced) after the last inserted paragraph —
# remove it so the paragraph count/structure matches the target exactly.
$lastIdx = $d.Paragraphs.Count
$last = $d.Paragraphs($lastIdx)
$lastText = $last.Range.Text
$lastStart = $last.Range.Start
$lastEnd = $last.Range.End
if ($lastText -eq "" -or $lastText -eq $null) {
    $delStart = $lastStart - 1
    $delRange = $d.Range($delStart, $lastEnd)
    $delRange.Delete()
}

Write-Output "done"
